$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Collapse the multi-run Title / Author / Abstract paragraphs into single
#    runs (the separate word-by-word runs get merged into one run each,
#    matching the author's "refresh" of the document front matter).
# ---------------------------------------------------------------------------

$titlePara = $d.Paragraphs.Item(1)
[void]$titlePara.Range.Find.Execute(
    "Questions: Multivariate implicit differentiation",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Questions: Multivariate implicit differentiation", 2)

$authorPara = $d.Paragraphs.Item(2)
[void]$authorPara.Range.Find.Execute(
    "Donald Campbell",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Donald Campbell", 2)

$abstractPara = $d.Paragraphs.Item(4)
[void]$abstractPara.Range.Find.Execute(
    "A selection of questions for the study guide on multivariate implicit differentiation.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A selection of questions for the study guide on multivariate implicit differentiation.", 2)

# ---------------------------------------------------------------------------
# 2. Normalise every <m:dPr> element so its children follow the canonical
#    schema order begChr, sepChr, endChr, grow (the original file had
#    endChr and sepChr swapped). We pull each paragraph's own OOXML via
#    WordOpenXML, fix up just the <m:dPr> child ordering with a regex, and
#    push the corrected paragraph back in with InsertXML (only on
#    paragraphs that actually contain a delimiter so everything else is
#    left completely untouched).
# ---------------------------------------------------------------------------

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $rng = $para.Range
    $xml = $rng.WordOpenXML

    if ($xml.Contains("<m:dPr>")) {
        $m = [regex]::Match($xml, '<w:p[ >].*?</w:p>')
        if ($m.Success) {
            $frag = $m.Value
            $fixedFrag = $frag -replace '<m:endChr([^/]*)/>\s*<m:sepChr([^/]*)/>', '<m:sepChr$2/><m:endChr$1/>'
            if ($fixedFrag -ne $frag) {
                [void]$rng.InsertXML($fixedFrag)
            }
        }
    }
}
